$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-70: column A previously held placeholder "1900-01-DD"/"19xx-.."/"20xx-.." style
# dates left over from a bad export. Restore the real period/category labels there.
# (A few of the restored labels look date-like to Excel's autosense, e.g. "2006  January",
# so those cells are pre-formatted as Text to stop them turning back into date serials.)
$ws.Range("A6:A10").NumberFormat = "@"

$newColumnA = [ordered]@{
    2 = '                  Central Bank of Sri Lanka'
    3 = '                  National Gem and Jewellery Authority'
    4 = '                  Sri Lanka Customs'
    5 = '  Sources: Ceylon Petroleum Corporation and Other Exporters of Petroleum'
    6 = '2006  January'
    7 = '2007  January'
    8 = '2008  January'
    9 = '2009   January'
    10 = '2010 January'
    11 = 'April'
    12 = 'April'
    13 = 'April'
    14 = 'April'
    15 = 'April'
    16 = 'Aug'
    17 = 'Aug'
    18 = 'Aug'
    19 = 'Aug'
    20 = 'Aug'
    21 = 'Dec'
    22 = 'Dec'
    23 = 'Dec'
    24 = 'Dec'
    25 = 'Dec'
    26 = 'February'
    27 = 'February'
    28 = 'February'
    29 = 'February'
    30 = 'February'
    31 = 'Industrial Exports'
    32 = 'July'
    33 = 'July'
    34 = 'July'
    35 = 'July'
    36 = 'July'
    37 = 'Jun'
    38 = 'Jun'
    39 = 'Jun'
    40 = 'Jun'
    41 = 'Jun'
    42 = 'March'
    43 = 'March'
    44 = 'March'
    45 = 'March'
    46 = 'March'
    47 = 'May'
    48 = 'May'
    49 = 'May'
    50 = 'May'
    51 = 'May'
    52 = 'Nov'
    53 = 'Nov'
    54 = 'Nov'
    55 = 'Nov'
    56 = 'Nov'
    57 = 'Oct'
    58 = 'Oct'
    59 = 'Oct'
    60 = 'Oct'
    61 = 'Oct'
    62 = 'Period'
    63 = 'Sep'
    64 = 'Sep'
    65 = 'Sep'
    66 = 'Sep'
    67 = 'Sep'
    68 = 'Table 2.02.5: Exports (US$ Million)'
    69 = 'Table 2.02: Exports - Monthly (2006-2010)'
    70 = 'Table 2.02: Exports - Monthly (2006-2010)'
}

foreach ($row in $newColumnA.Keys) {
    $ws.Cells.Item($row, 1).Value = $newColumnA[$row]
}

# Rows 71-224 held a stray alphabetically-sorted dump of every label/date used on the sheet
# (an artifact of the bad export) - clear it all out so only the genuine 2-70 rows remain.
$ws.Range("A71:S224").ClearContents()

# Keep rows 71-98 as part of the sheet's used range (matching the trimmed sheet dimensions)
# by touching column A there, then clearing the formatting straight back to the default.
$ws.Range("A71:A98").Font.Name = "Calibri"
$ws.Range("A71:A98").ClearFormats()
